$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# 1. Update C26 value (1.1 -> -0.08)
$ws.Range("C26").Value = -0.08

# 2. Insert 7 new rows before the old row 28 ("Main gear" / x_mg), pushing it to row 35
$ws.Range("A28:A34").EntireRow.Insert()

# 3. Fill in the new rows 28-34
$ws.Range("B28").Value = "rootChordLen_w"
$ws.Range("C28").Value = 8.1
$ws.Range("D28").Value = "m"
$ws.Range("E28").Value = "Root chord Length"

$ws.Range("B29").Value = "tipChordLen_w"
$ws.Range("C29").Value = 1.62
$ws.Range("D29").Value = "m"
$ws.Range("E29").Value = "Tip chord Length"

$ws.Range("B30").Value = "wingSemiSpan"
$ws.Range("C30").Formula = "=36.41/2"
$ws.Range("D30").Value = "m"
$ws.Range("E30").Value = "Semi span of main wing"

$ws.Range("B31").Value = "fuelTankLen"
$ws.Range("C31").Value = 0.8
$ws.Range("D31").Value = "-"
$ws.Range("E31").Value = "Pct of main wing semi-span"

$ws.Range("B32").Value = "engineLoc_1"
$ws.Range("C32").Value = 7.075
$ws.Range("D32").Value = "m"
$ws.Range("E32").Value = "Length of engine from root of wing"

$ws.Range("B33").Value = "engineLoc_2"
$ws.Range("C33").Value = 11.875
$ws.Range("D33").Value = "m"
$ws.Range("E33").Value = "Length of engine from root of wing"

$ws.Range("B34").Value = "engineWeight"
$ws.Range("C34").Value = 2177
$ws.Range("D34").Value = "kg"
$ws.Range("E34").Value = "Weight of engine"

# 4. Apply number format style (matching column C's usual style) to the new C cells,
#    and the "applyFill" marker style on C31 specifically.
$ws.Range("C28").Copy()
$ws.Range("C29:C30").PasteSpecial(-4122)
$ws.Range("C32:C34").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C31").Interior.ColorIndex = -4142

# 5. Update selection / view on the Data sheet
$ws.Application.ActiveWindow.ScrollRow = 11
$ws.Range("D29").Select()
